# Applies the cryptos-list price/volume refresh described by the commit:
# "Updated cryptos list ... with GitHub Actions" -- column D (Price) and
# column E (Volume(1h)) on rows 2-51 get new scraped text values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) ---
# Values that cannot be misread as a plain number (multiple '.' separators,
# or containing the subscript-3 character used for PEPE's price) can be
# written straight through; Excel keeps them as literal text.
$ws.Cells.Item(2, 4).Value2 = '56.874.28'
$ws.Cells.Item(3, 4).Value2 = '2.976.56'
$ws.Cells.Item(12, 4).Value2 = '3.496.60'
$ws.Cells.Item(16, 4).Value2 = '56.986.97'
$ws.Cells.Item(18, 4).Value2 = '2.975.94'
$ws.Cells.Item(28, 4).Value2 = '0.0₃0892'
$ws.Cells.Item(40, 4).Value2 = '3.008.71'
$ws.Cells.Item(45, 4).Value2 = '2.196.10'

# Values that DO look like a plain number ("1.00", "0.428", ...) need the
# cell pre-formatted as Text, otherwise Excel's type inference silently
# converts them to numeric values and the exact literal digits (trailing
# zeros etc.) would not round-trip. Style is restored to Normal afterwards
# so the cell keeps its original (unstyled) appearance.
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.00'
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value2 = '498.30'
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value2 = '137.25'
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.428'
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value2 = '7.40'
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.357'
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value2 = '25.82'
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value2 = '6.08'
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value2 = '12.62'
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value2 = '7.78'
$c.Style = "Normal"
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value2 = '320.94'
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.997'
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.487'
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value2 = '63.57'
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value2 = '6.56'
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value2 = '7.12'
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value2 = '20.18'
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value2 = '153.18'
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value2 = '4.64'
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value2 = '5.75'
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.24'
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value2 = '23.95'
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value2 = '37.44'
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.00'
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value2 = '3.74'
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.640'
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.38'
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.948'
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value2 = '5.94'
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value2 = '19.11'
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.80'
$c.Style = "Normal"

# --- Volume(1h) (column E) ---
# Every value carries '%' / sign / padding spaces, so none of these are
# ever misread as numbers.
$ws.Cells.Item(2, 5).Value2 = '  +0.12%  '
$ws.Cells.Item(3, 5).Value2 = '  -0.90%  '
$ws.Cells.Item(4, 5).Value2 = '  +0.11%  '
$ws.Cells.Item(5, 5).Value2 = '  -2.60%  '
$ws.Cells.Item(6, 5).Value2 = '  -1.52%  '
$ws.Cells.Item(7, 5).Value2 = '  +0.26%  '
$ws.Cells.Item(8, 5).Value2 = '  -1.77%  '
$ws.Cells.Item(9, 5).Value2 = '  -0.94%  '
$ws.Cells.Item(10, 5).Value2 = '  -0.30%  '
$ws.Cells.Item(11, 5).Value2 = '  -0.17%  '
$ws.Cells.Item(12, 5).Value2 = '  -0.62%  '
$ws.Cells.Item(13, 5).Value2 = '  -1.17%  '
$ws.Cells.Item(14, 5).Value2 = '  -0.63%  '
$ws.Cells.Item(15, 5).Value2 = '  +1.09%  '
$ws.Cells.Item(16, 5).Value2 = '  +0.27%  '
$ws.Cells.Item(17, 5).Value2 = '  +2.28%  '
$ws.Cells.Item(18, 5).Value2 = '  -0.91%  '
$ws.Cells.Item(19, 5).Value2 = '  +0.56%  '
$ws.Cells.Item(20, 5).Value2 = '  -0.95%  '
$ws.Cells.Item(21, 5).Value2 = '  -1.86%  '
$ws.Cells.Item(22, 5).Value2 = '  -0.26%  '
$ws.Cells.Item(23, 5).Value2 = '  -0.93%  '
$ws.Cells.Item(24, 5).Value2 = '  -0.07%  '
$ws.Cells.Item(25, 5).Value2 = '  +0.17%  '
$ws.Cells.Item(26, 5).Value2 = '  +0.91%  '
$ws.Cells.Item(27, 5).Value2 = '  -4.74%  '
$ws.Cells.Item(28, 5).Value2 = '  -2.15%  '
$ws.Cells.Item(29, 5).Value2 = '  -1.70%  '
$ws.Cells.Item(30, 5).Value2 = '  +0.93%  '
$ws.Cells.Item(31, 5).Value2 = '  -2.29%  '
$ws.Cells.Item(32, 5).Value2 = '  -5.96%  '
$ws.Cells.Item(33, 5).Value2 = '  -1.96%  '
$ws.Cells.Item(34, 5).Value2 = '  -1.13%  '
$ws.Cells.Item(35, 5).Value2 = '  +1.32%  '
$ws.Cells.Item(36, 5).Value2 = '  +0.65%  '
$ws.Cells.Item(37, 5).Value2 = '  -2.47%  '
$ws.Cells.Item(38, 5).Value2 = '  -0.06%  '
$ws.Cells.Item(39, 5).Value2 = '  -2.06%  '
$ws.Cells.Item(40, 5).Value2 = '  -0.90%  '
$ws.Cells.Item(41, 5).Value2 = '  +0.88%  '
$ws.Cells.Item(42, 5).Value2 = '  +0.13%  '
$ws.Cells.Item(43, 5).Value2 = '  +1.56%  '
$ws.Cells.Item(44, 5).Value2 = '  -1.27%  '
$ws.Cells.Item(45, 5).Value2 = '  -4.28%  '
$ws.Cells.Item(46, 5).Value2 = '  -2.65%  '
$ws.Cells.Item(47, 5).Value2 = '  -5.70%  '
$ws.Cells.Item(48, 5).Value2 = '  +0.81%  '
$ws.Cells.Item(49, 5).Value2 = '  -1.80%  '
$ws.Cells.Item(50, 5).Value2 = '  -1.48%  '
$ws.Cells.Item(51, 5).Value2 = '  -7.51%  '
